$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(5, 6).Value = 3041
$ws1.Cells.Item(9, 6).Value = 300
$ws1.Cells.Item(10, 6).Value = 7159
$ws1.Cells.Item(11, 6).Value = 54
$ws1.Cells.Item(14, 6).Value = 642
$ws1.Cells.Item(15, 6).Value = 1551
$ws1.Cells.Item(16, 6).Value = 2312
$ws1.Cells.Item(18, 6).Value = 1161
$ws1.Cells.Item(19, 6).Value = 29
$ws1.Cells.Item(24, 6).Value = 1846
$ws1.Cells.Item(25, 6).Value = 1742
$ws1.Cells.Item(28, 6).Value = 1706
$ws1.Cells.Item(29, 6).Value = 1294
$ws1.Cells.Item(31, 6).Value = 600
$ws1.Cells.Item(32, 6).Value = 26
$ws1.Cells.Item(33, 6).Value = 1081
$ws1.Cells.Item(35, 6).Value = 45
$ws1.Cells.Item(37, 6).Value = 2858
$ws1.Cells.Item(38, 6).Value = 2109
$ws1.Cells.Item(39, 6).Value = 65
$ws1.Cells.Item(40, 6).Value = 200
$ws1.Cells.Item(43, 6).Value = 26
$ws1.Cells.Item(44, 6).Value = 37
$ws1.Cells.Item(46, 6).Value = 135
$ws1.Cells.Item(47, 6).Value = 202
$ws1.Cells.Item(48, 6).Value = 175
$ws1.Cells.Item(49, 6).Value = 73

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(10, 6).Value = 67
$ws2.Cells.Item(18, 6).Value = 344
$ws2.Cells.Item(19, 6).Value = 501
$ws2.Cells.Item(25, 6).Value = 20
$ws2.Cells.Item(27, 6).Value = 28

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(4, 6).Value = 553
$ws3.Cells.Item(6, 6).Value = 1765
$ws3.Cells.Item(8, 6).Value = 2828
$ws3.Cells.Item(9, 6).Value = 1079
$ws3.Cells.Item(10, 6).Value = 1018
$ws3.Cells.Item(12, 6).Value = 365
$ws3.Cells.Item(13, 6).Value = 1702
$ws3.Cells.Item(14, 6).Value = 7816

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(4, 6).Value = 3041
$ws4.Cells.Item(6, 6).Value = 1765
$ws4.Cells.Item(7, 6).Value = 300
$ws4.Cells.Item(8, 6).Value = 2828
$ws4.Cells.Item(9, 6).Value = 7159
$ws4.Cells.Item(10, 6).Value = 1079
$ws4.Cells.Item(11, 6).Value = 54
$ws4.Cells.Item(13, 6).Value = 365
$ws4.Cells.Item(14, 6).Value = 642
$ws4.Cells.Item(15, 6).Value = 1551
$ws4.Cells.Item(16, 6).Value = 2312
$ws4.Cells.Item(18, 6).Value = 1161
$ws4.Cells.Item(19, 6).Value = 29
$ws4.Cells.Item(21, 6).Value = 67
$ws4.Cells.Item(23, 6).Value = 1846
$ws4.Cells.Item(27, 6).Value = 1706
$ws4.Cells.Item(28, 6).Value = 1294
$ws4.Cells.Item(30, 6).Value = 600
$ws4.Cells.Item(31, 6).Value = 26
$ws4.Cells.Item(32, 6).Value = 1081
$ws4.Cells.Item(35, 6).Value = 501
$ws4.Cells.Item(37, 6).Value = 45
$ws4.Cells.Item(39, 6).Value = 2858
$ws4.Cells.Item(40, 6).Value = 2109
$ws4.Cells.Item(41, 6).Value = 65
$ws4.Cells.Item(42, 6).Value = 200
$ws4.Cells.Item(46, 6).Value = 20
$ws4.Cells.Item(47, 6).Value = 202
